$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("M2").Value = 1.08
$ws.Range("O2").Value = 1.47
$ws.Range("S2").Value = 1.54
$ws.Range("BD2").Value = 151

# Row 3 updates
$ws.Range("G3").Value = 1.7
$ws.Range("L3").Value = 6
$ws.Range("AJ3").Value = 67
$ws.Range("AR3").Value = 51
$ws.Range("AV3").Value = 81
$ws.Range("BD3").Value = 126
